$wb = $excel.ActiveWorkbook

# ALC row 33 (Leve Item ID 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 18921.062
$ws.Cells.Item(33, 9).Value = 22360.424
$ws.Cells.Item(33, 11).Value = 22360.424
$ws.Cells.Item(33, 13).Value = -22131.424

# ALC row 38 (Leve Item ID 4599)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(38, 8).Value = 1673.375
$ws.Cells.Item(38, 10).Value = 2659.2
$ws.Cells.Item(38, 12).Value = 7977.599999999999
$ws.Cells.Item(38, 14).Value = -8721.599999999999

# ALC row 39 (Leve Item ID 4603)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(39, 8).Value = 14525.5
$ws.Cells.Item(39, 9).Value = 1568.3334
$ws.Cells.Item(39, 11).Value = 4705.0002
$ws.Cells.Item(39, 13).Value = -4409.0002

# ALC row 40 (Leve Item ID 5505)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 1530.2
$ws.Cells.Item(40, 9).Value = 1679.2
$ws.Cells.Item(40, 10).Value = 1480.5333
$ws.Cells.Item(40, 11).Value = 1679.2
$ws.Cells.Item(40, 12).Value = 1480.5333
$ws.Cells.Item(40, 13).Value = -1504.2
$ws.Cells.Item(40, 14).Value = -1830.5333

# ALC row 43 (Leve Item ID 5472)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 9102253
$ws.Cells.Item(43, 9).Value = 25002500
$ws.Cells.Item(43, 10).Value = 16398.715
$ws.Cells.Item(43, 11).Value = 25002500
$ws.Cells.Item(43, 12).Value = 16398.715
$ws.Cells.Item(43, 13).Value = -25002431
$ws.Cells.Item(43, 14).Value = -16536.715

# ALC row 137 (Leve Item ID 44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3854399
$ws.Cells.Item(137, 9).Value = 6132
$ws.Cells.Item(137, 11).Value = 18396
$ws.Cells.Item(137, 13).Value = -15846

# ARM row 45 (Leve Item ID 27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 93462.27
$ws.Cells.Item(45, 9).Value = 168847.5
$ws.Cells.Item(45, 11).Value = 168847.5
$ws.Cells.Item(45, 13).Value = -168470.5

# ARM row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3049005
$ws.Cells.Item(61, 9).Value = 71348.44
$ws.Cells.Item(61, 11).Value = 71348.44
$ws.Cells.Item(61, 13).Value = -71136.44

# ARM row 102 (Leve Item ID 19945)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(102, 8).Value = 1975.0952
$ws.Cells.Item(102, 9).Value = 1975.0952
$ws.Cells.Item(102, 11).Value = 1975.0952
$ws.Cells.Item(102, 13).Value = -353.0952

# ARM row 110 (Leve Item ID 27708)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 695.4358999999999
$ws.Cells.Item(110, 9).Value = 686.5278
$ws.Cells.Item(110, 11).Value = 686.5278
$ws.Cells.Item(110, 13).Value = 1358.4722

# ARM row 119 (Leve Item ID 26287)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(119, 8).Value = 69333
$ws.Cells.Item(119, 10).Value = 69333
$ws.Cells.Item(119, 12).Value = 69333
$ws.Cells.Item(119, 14).Value = -79009

# ARM row 122 (Leve Item ID 36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2085.5881
$ws.Cells.Item(122, 9).Value = 1903.4375
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 5710.3125
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -3260.3125
$ws.Cells.Item(122, 14).Value = -19900

# ARM row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3049005
$ws.Cells.Item(136, 9).Value = 71348.44
$ws.Cells.Item(136, 11).Value = 214045.32
$ws.Cells.Item(136, 13).Value = -211495.32

# BSM row 20 (Leve Item ID 14149)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 852.1177
$ws.Cells.Item(20, 9).Value = 824.0909
$ws.Cells.Item(20, 10).Value = 903.5
$ws.Cells.Item(20, 11).Value = 824.0909
$ws.Cells.Item(20, 12).Value = 903.5
$ws.Cells.Item(20, 13).Value = -577.0909
$ws.Cells.Item(20, 14).Value = -1397.5

# BSM row 99 (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 11587.167
$ws.Cells.Item(99, 10).Value = 2709.5
$ws.Cells.Item(99, 12).Value = 2709.5
$ws.Cells.Item(99, 14).Value = -5705.5

# CRP row 9 (Leve Item ID 15611)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(9, 8).Value = 31999.5
$ws.Cells.Item(9, 10).Value = 31999.5
$ws.Cells.Item(9, 12).Value = 31999.5
$ws.Cells.Item(9, 14).Value = -32335.5

# CRP row 14 (Leve Item ID 1998)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(14, 8).Value = 6450
$ws.Cells.Item(14, 9).Value = 4900
$ws.Cells.Item(14, 10).Value = 8000
$ws.Cells.Item(14, 11).Value = 4900
$ws.Cells.Item(14, 12).Value = 8000
$ws.Cells.Item(14, 13).Value = -4730
$ws.Cells.Item(14, 14).Value = -8340

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2456.9363
$ws.Cells.Item(31, 9).Value = 2810.762
$ws.Cells.Item(31, 11).Value = 2810.762
$ws.Cells.Item(31, 13).Value = -2515.762

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(34, 8).Value = 2456.9363
$ws.Cells.Item(34, 9).Value = 2810.762
$ws.Cells.Item(34, 11).Value = 2810.762
$ws.Cells.Item(34, 13).Value = -2608.762

# CUL row 5 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2354.7778
$ws.Cells.Item(5, 9).Value = 1923.25
$ws.Cells.Item(5, 11).Value = 5769.75
$ws.Cells.Item(5, 13).Value = -5657.75

# CUL row 17 (Leve Item ID 4640)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 1850
$ws.Cells.Item(17, 10).Value = 3199.5
$ws.Cells.Item(17, 12).Value = 9598.5
$ws.Cells.Item(17, 14).Value = -9936.5

# CUL row 46 (Leve Item ID 4701)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 100002340
$ws.Cells.Item(46, 10).Value = 250003700
$ws.Cells.Item(46, 12).Value = 750011100
$ws.Cells.Item(46, 14).Value = -750011282

# CUL row 56 (Leve Item ID 10146)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 12825913
$ws.Cells.Item(56, 9).Value = 12825913
$ws.Cells.Item(56, 11).Value = 12825913
$ws.Cells.Item(56, 13).Value = -12825383

# CUL row 87 (Leve Item ID 12864)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 14942
$ws.Cells.Item(87, 9).Value = 2437
$ws.Cells.Item(87, 10).Value = 24946
$ws.Cells.Item(87, 11).Value = 7311
$ws.Cells.Item(87, 12).Value = 74838
$ws.Cells.Item(87, 13).Value = -6063
$ws.Cells.Item(87, 14).Value = -77334

# CUL row 90 (Leve Item ID 12864)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(90, 8).Value = 14942
$ws.Cells.Item(90, 9).Value = 2437
$ws.Cells.Item(90, 10).Value = 24946
$ws.Cells.Item(90, 11).Value = 21933
$ws.Cells.Item(90, 12).Value = 224514
$ws.Cells.Item(90, 13).Value = -15693
$ws.Cells.Item(90, 14).Value = -236994

# CUL row 135 (Leve Item ID 43974)
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(135, 8).Value = 2354.7778
$ws.Cells.Item(135, 9).Value = 1923.25
$ws.Cells.Item(135, 11).Value = 17309.25
$ws.Cells.Item(135, 13).Value = -14774.25

# GSM row 2 (Leve Item ID 5062)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 6321.75
$ws.Cells.Item(2, 9).Value = 64.333336
$ws.Cells.Item(2, 11).Value = 64.333336
$ws.Cells.Item(2, 13).Value = 48.666664

# GSM row 70 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4750.1665
$ws.Cells.Item(70, 9).Value = 4811.75
$ws.Cells.Item(70, 10).Value = 4627
$ws.Cells.Item(70, 11).Value = 4811.75
$ws.Cells.Item(70, 12).Value = 4627
$ws.Cells.Item(70, 13).Value = -4541.75
$ws.Cells.Item(70, 14).Value = -5167

# GSM row 73 (Leve Item ID 14146)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(73, 8).Value = 4750.1665
$ws.Cells.Item(73, 9).Value = 4811.75
$ws.Cells.Item(73, 10).Value = 4627
$ws.Cells.Item(73, 11).Value = 4811.75
$ws.Cells.Item(73, 12).Value = 4627
$ws.Cells.Item(73, 13).Value = -3875.75
$ws.Cells.Item(73, 14).Value = -6499

# GSM row 122 (Leve Item ID 36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1843
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).ClearContents()

# GSM row 125 (Leve Item ID 34291)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(125, 8).Value = 75000
$ws.Cells.Item(125, 10).Value = 75000
$ws.Cells.Item(125, 12).Value = 75000
$ws.Cells.Item(125, 14).Value = -79920

# GSM row 126 (Leve Item ID 36184)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3296.6924
$ws.Cells.Item(126, 9).Value = 3323.5454
$ws.Cells.Item(126, 10).Value = 3149
$ws.Cells.Item(126, 11).Value = 9970.636200000001
$ws.Cells.Item(126, 12).Value = 9447
$ws.Cells.Item(126, 13).Value = -7500.636200000001
$ws.Cells.Item(126, 14).Value = -14387

# GSM row 139 (Leve Item ID 42373)
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(139, 8).Value = 206570
$ws.Cells.Item(139, 10).Value = 251399
$ws.Cells.Item(139, 12).Value = 251399
$ws.Cells.Item(139, 14).Value = -261679

# LTW row 36 (Leve Item ID 34261)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(36, 8).Value = 75000
$ws.Cells.Item(36, 10).Value = 75000
$ws.Cells.Item(36, 12).Value = 75000
$ws.Cells.Item(36, 14).Value = -76124

# LTW row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1166.2727
$ws.Cells.Item(55, 9).Value = 1342.0625
$ws.Cells.Item(55, 10).Value = 1000.82355
$ws.Cells.Item(55, 11).Value = 1342.0625
$ws.Cells.Item(55, 12).Value = 1000.82355
$ws.Cells.Item(55, 13).Value = -1169.0625
$ws.Cells.Item(55, 14).Value = -1346.82355

# LTW row 80 (Leve Item ID 12027)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

# LTW row 83 (Leve Item ID 12027)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

# LTW row 124 (Leve Item ID 34264)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(124, 8).Value = 75000
$ws.Cells.Item(124, 10).Value = 75000
$ws.Cells.Item(124, 12).Value = 75000
$ws.Cells.Item(124, 14).Value = -84820

# LTW row 132 (Leve Item ID 44058)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 2851.6296
$ws.Cells.Item(132, 9).Value = 2494.4119
$ws.Cells.Item(132, 11).Value = 7483.2357
$ws.Cells.Item(132, 13).Value = -4953.2357

# LTW row 134 (Leve Item ID 42024)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(134, 8).Value = 75000
$ws.Cells.Item(134, 10).Value = 75000
$ws.Cells.Item(134, 12).Value = 75000
$ws.Cells.Item(134, 14).Value = -85140

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 74536.64
$ws.Cells.Item(136, 10).Value = 4946.5
$ws.Cells.Item(136, 12).Value = 14839.5
$ws.Cells.Item(136, 14).Value = -19939.5

# LTW row 140 (Leve Item ID 42503)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(140, 8).Value = 124997.5
$ws.Cells.Item(140, 10).Value = 124997.5
$ws.Cells.Item(140, 12).Value = 124997.5
$ws.Cells.Item(140, 14).Value = -135357.5

# LTW row 141 (Leve Item ID 42487)
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(141, 8).Value = 169995
$ws.Cells.Item(141, 10).Value = 169995
$ws.Cells.Item(141, 12).Value = 169995
$ws.Cells.Item(141, 14).Value = -180355

# WVR row 81 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 128416.625
$ws.Cells.Item(81, 9).Value = 4169.2
$ws.Cells.Item(81, 11).Value = 8338.4
$ws.Cells.Item(81, 13).Value = -7277.4

# WVR row 84 (Leve Item ID 12596)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(84, 8).Value = 128416.625
$ws.Cells.Item(84, 9).Value = 4169.2
$ws.Cells.Item(84, 11).Value = 41692
$ws.Cells.Item(84, 13).Value = -36388

# WVR row 100 (Leve Item ID 19981)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 1097.3636
$ws.Cells.Item(100, 9).Value = 896.7778
$ws.Cells.Item(100, 11).Value = 1793.5556
$ws.Cells.Item(100, 13).Value = -1252.5556

# WVR row 128 (Leve Item ID 34563)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(128, 8).Value = 95994.39999999999
$ws.Cells.Item(128, 10).Value = 95994.39999999999
$ws.Cells.Item(128, 12).Value = 95994.39999999999
$ws.Cells.Item(128, 14).Value = -105954.4

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 2461.8
$ws.Cells.Item(132, 9).Value = 1733.8334
$ws.Cells.Item(132, 11).Value = 5201.5002
$ws.Cells.Item(132, 13).Value = -2671.5002

Write-Output "Applied all Balmung_Profits updates"